$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Add an "Email" column to the PROGRAMMERS table (columns G:I) and to the
# USERS table (columns N:P), each with mailto: hyperlinks to the matching
# programmer / user.
# ---------------------------------------------------------------------------

# --- Headers -----------------------------------------------------------
$ws.Range("I11").Value = "Email"
$ws.Range("H11").Copy() | Out-Null
$ws.Range("I11").PasteSpecial(-4122) | Out-Null

$ws.Range("P11").Value = "Email"
$ws.Range("O11").Copy() | Out-Null
$ws.Range("P11").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# --- Prime borders for the new cells -----------------------------------
# I12 / P12 need a full thin box border, matching the existing boxed cells
# used elsewhere in the data tables (e.g. E13).
$ws.Range("E13").Copy() | Out-Null
$ws.Range("I12").PasteSpecial(-4122) | Out-Null
$ws.Range("P12").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# I13:I15 and P13 continue the box without a bottom edge.
$ws.Range("I12").Copy() | Out-Null
$ws.Range("I13:I15").PasteSpecial(-4122) | Out-Null
$ws.Range("P13").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0
# NOTE: assigning Borders.Item(..).LineStyle on a multi-cell range only takes
# effect on the last cell of the range in this engine, so each cell needs to
# be addressed individually.
foreach ($addr in @("I13", "I14", "I15", "P13")) {
  $ws.Range($addr).Borders.Item(9).LineStyle = -4142
}

# P14 / P15 only keep a thin top edge (bottom of the USERS table).
foreach ($addr in @("P14", "P15")) {
  $ws.Range($addr).Borders.Item(8).LineStyle = 1
  $ws.Range($addr).Borders.Item(8).Weight = 2
}

# --- Hyperlinks / values -------------------------------------------------
$ws.Hyperlinks.Add($ws.Range("I12"), "mailto:Luke@l.com", "", "", "Luke@l.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I13"), "mailto:Mitch@m.com", "", "", "Mitch@m.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I14"), "mailto:Ryan@r.com", "", "", "Ryan@r.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("I15"), "mailto:Gearoid@g.com", "", "", "Gearoid@g.com") | Out-Null

$ws.Hyperlinks.Add($ws.Range("P12"), "mailto:Luke@l.com", "", "", "Luke@l.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("P13"), "mailto:Mitch@m.com", "", "", "Mitch@m.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("P14"), "mailto:Ryan@r.com", "", "", "Ryan@r.com") | Out-Null

# P15 stays blank (USERS table has only 3 data rows).

# --- Column width for the new column -----------------------------------
$ws.Columns.Item(16).ColumnWidth = 14

# --- View state ----------------------------------------------------------
$ws.Range("R17").Select()

Write-Host "Email columns with hyperlinks added"
